$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  6/30/2025  Through  7/6/2025"

# --- Updated weekly crime-statistics figures ---
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = -50
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -78.571428571428
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 200
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 26
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = 44.444444444444
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 188.888888888889
$ws.Range("N15").Value = -33.333333333333
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = -71.428571428571
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 50
$ws.Range("H16").Value = -48
$ws.Range("I16").Value = 164
$ws.Range("J16").Value = 219
$ws.Range("K16").Value = -25.114155251141
$ws.Range("L16").Value = -4.093567251461
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -80.240963855421
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -21.428571428571
$ws.Range("F17").Value = 53
$ws.Range("H17").Value = -29.333333333333
$ws.Range("I17").Value = 380
$ws.Range("J17").Value = 386
$ws.Range("K17").Value = -1.554404145077
$ws.Range("L17").Value = 6.741573033707
$ws.Range("M17").Value = 126.190476190476
$ws.Range("N17").Value = -13.043478260869
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 30.76923076923
$ws.Range("I18").Value = 117
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 30
$ws.Range("L18").Value = 40.963855421686
$ws.Range("M18").Value = -8.59375
$ws.Range("N18").Value = -80.102040816326
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 30
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 2.173913043478
$ws.Range("I19").Value = 269
$ws.Range("J19").Value = 324
$ws.Range("K19").Value = -16.975308641975
$ws.Range("L19").Value = -3.237410071942
$ws.Range("M19").Value = 42.328042328042
$ws.Range("N19").Value = -54.789915966386
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 155.555555555556
$ws.Range("I20").Value = 106
$ws.Range("J20").Value = 99
$ws.Range("K20").Value = 7.070707070707
$ws.Range("L20").Value = -12.396694214876
$ws.Range("M20").Value = 15.217391304347
$ws.Range("N20").Value = -85.828877005347
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -4.761904761904
$ws.Range("F21").Value = 175
$ws.Range("G21").Value = 199
$ws.Range("H21").Value = -12.060301507537
$ws.Range("I21").Value = 1065
$ws.Range("J21").Value = 1139
$ws.Range("K21").Value = -6.49692712906
$ws.Range("L21").Value = 4.207436399217
$ws.Range("M21").Value = 27.08830548926
$ws.Range("N21").Value = -67.240848969547
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = 11.764705882352
$ws.Range("M22").Value = 58.333333333333
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 31
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = 19.230769230769
$ws.Range("L23").Value = 47.619047619047
$ws.Range("M23").Value = 47.619047619047
$ws.Range("C24").Value = 42
$ws.Range("D24").Value = 55
$ws.Range("E24").Value = -23.636363636363
$ws.Range("F24").Value = 155
$ws.Range("G24").Value = 199
$ws.Range("H24").Value = -22.110552763819
$ws.Range("I24").Value = 994
$ws.Range("J24").Value = 960
$ws.Range("K24").Value = 3.541666666666
$ws.Range("L24").Value = 22.413793103448
$ws.Range("M24").Value = 64.842454394693
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 45
$ws.Range("E25").Value = -35.555555555555
$ws.Range("F25").Value = 90
$ws.Range("G25").Value = 147
$ws.Range("H25").Value = -38.775510204081
$ws.Range("I25").Value = 569
$ws.Range("J25").Value = 611
$ws.Range("K25").Value = -6.873977086743
$ws.Range("L25").Value = 68.343195266272
$ws.Range("C26").Value = 29
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = 45
$ws.Range("F26").Value = 98
$ws.Range("G26").Value = 88
$ws.Range("H26").Value = 11.363636363636
$ws.Range("I26").Value = 478
$ws.Range("J26").Value = 525
$ws.Range("K26").Value = -8.95238095238
$ws.Range("L26").Value = -5.719921104536
$ws.Range("M26").Value = 14.62829736211
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 60
$ws.Range("I27").Value = 27
$ws.Range("J27").Value = 26
$ws.Range("K27").Value = 3.846153846153
$ws.Range("L27").Value = 28.571428571428
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 20
$ws.Range("H28").Value = -45
$ws.Range("I28").Value = 63
$ws.Range("J28").Value = 69
$ws.Range("K28").Value = -8.695652173913
$ws.Range("L28").Value = 40
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("M29").Value = -48
$ws.Range("N29").Value = -83.116883116883
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("M30").Value = -75
$ws.Range("N30").Value = -94.444444444444
$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 3
$ws.Range("L31").Value = -40
